$d = $word.ActiveDocument

$pairs = @(
    @("549×3=", "995×2="),
    @("849×9=", "243×5="),
    @("285×9=", "646×8="),
    @("456×8=", "768×8="),
    @("561×4=", "680×9="),
    @("670×5=", "556×7="),
    @("800×8=", "433×5="),
    @("283×6=", "612×2="),
    @("945×3=", "241×8="),
    @("880×6=", "957×7="),
    @("551×9=", "730×6="),
    @("221×6=", "576×6="),
    @("715×8=", "248×7="),
    @("152×8=", "907×8="),
    @("655×7=", "228×4="),
    @("106×3=", "747×6="),
    @("509×7=", "913×7="),
    @("551×2=", "320×6="),
    @("527×9=", "866×3="),
    @("361×6=", "935×7="),
    @("584×2=", "212×3="),
    @("396×7=", "139×6="),
    @("232×2=", "641×5="),
    @("884×2=", "870×4="),
    @("526×6=", "319×4=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
